$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.358186072990634
$ws.Range("C2").Value = 0.03435164755344999
$ws.Range("D2").Value = 0.4825267819929024
$ws.Range("E2").Value = 0.1405418174164286
$ws.Range("G2").Value = 2.395263194951099
$ws.Range("H2").Value = 1.966873903713463
$ws.Range("I2").Value = 1.937923625467633
$ws.Range("J2").Value = 0.04188490355932917
$ws.Range("K2").Value = 1.036683649955421
$ws.Range("L2").Value = 0.536819484399345
$ws.Range("N2").Value = 3.11950807204979

$ws.Range("B3").Value = 1.322521896112931
$ws.Range("C3").Value = 0.03000589728560499
$ws.Range("D3").Value = 0.4802977967310227
$ws.Range("E3").Value = 0.1406360495989247
$ws.Range("G3").Value = 2.396181371587545
$ws.Range("H3").Value = 1.972882305250238
$ws.Range("I3").Value = 1.943688605464473
$ws.Range("J3").Value = 0.04197150157384888
$ws.Range("K3").Value = 0.9994991532403219
$ws.Range("L3").Value = 0.5328514373847639
$ws.Range("N3").Value = 3.141374819484959

$ws.Range("B4").Value = 1.301299087843944
$ws.Range("C4").Value = 0.02733937641094997
$ws.Range("D4").Value = 0.4791285239344631
$ws.Range("E4").Value = 0.1407397137084914
$ws.Range("G4").Value = 2.397784435433266
$ws.Range("H4").Value = 1.977256265914662
$ws.Range("I4").Value = 1.947948332756972
$ws.Range("J4").Value = 0.0420280558851478
$ws.Range("K4").Value = 0.9772068731510899
$ws.Range("L4").Value = 0.5306493398677787
$ws.Range("N4").Value = 3.155597615098401

$ws.Range("B5").Value = 1.292820888855829
$ws.Range("C5").Value = 0.02625316874153327
$ws.Range("D5").Value = 0.4787022508438383
$ws.Range("E5").Value = 0.1407935008334622
$ws.Range("G5").Value = 2.398699065556059
$ws.Range("H5").Value = 1.979211031126383
$ws.Range("I5").Value = 1.949865351541881
$ws.Range("J5").Value = 0.04205195503304981
$ws.Range("K5").Value = 0.9682585057409199
$ws.Range("L5").Value = 0.5298109853888633
$ws.Range("N5").Value = 3.161593829268995

$ws.Range("B6").Value = 1.291423388869475
$ws.Range("C6").Value = 0.02607282940918765
$ws.Range("D6").Value = 0.4786345040824358
$ws.Range("E6").Value = 0.1408031300627535
$ws.Range("G6").Value = 2.398866726128276
$ws.Range("H6").Value = 1.979546031742458
$ws.Range("I6").Value = 1.950194614590224
$ws.Range("J6").Value = 0.04205597503987146
$ws.Range("K6").Value = 0.9667808554555108
$ws.Range("L6").Value = 0.5296753448277514
$ws.Range("N6").Value = 3.16260159481979

$ws.Range("B7").Value = 1.301184057818773
$ws.Range("C7").Value = 0.02732472576529688
$ws.Range("D7").Value = 0.4791225716263483
$ws.Range("E7").Value = 0.1407403923290662
$ws.Range("G7").Value = 2.397795712150241
$ws.Range("H7").Value = 1.977281930570499
$ws.Range("I7").Value = 1.94797345277199
$ws.Range("J7").Value = 0.04202837474177645
$ws.Range("K7").Value = 0.9770856415562434
$ws.Range("L7").Value = 0.5306377944276335
$ws.Range("N7").Value = 3.155677671079744

$ws.Range("B8").Value = 1.345749222942999
$ws.Range("C8").Value = 0.0328528436194091
$ws.Range("D8").Value = 0.481716909046142
$ws.Range("E8").Value = 0.1405648165970845
$ws.Range("G8").Value = 2.395364086511378
$ws.Range("H8").Value = 1.968803561121348
$ws.Range("I8").Value = 1.939762006928802
$ws.Range("J8").Value = 0.04191406176931034
$ws.Range("K8").Value = 1.023750726033256
$ws.Range("L8").Value = 0.5354027312681353
$ws.Range("N8").Value = 3.126882419420852

$ws.Range("B9").Value = 1.43848255459892
$ws.Range("C9").Value = 0.04370945098318657
$ws.Range("D9").Value = 0.4883827840613009
$ws.Range("E9").Value = 0.1405829555438558
$ws.Range("G9").Value = 2.39884341644192
$ws.Range("H9").Value = 1.957605429073027
$ws.Range("I9").Value = 1.929369330203251
$ws.Range("J9").Value = 0.04171663718693219
$ws.Range("K9").Value = 1.11952963905452
$ws.Range("L9").Value = 0.5466024350472196
$ws.Range("N9").Value = 3.076733767715929

$ws.Range("B10").Value = 1.509857735194998
$ws.Range("C10").Value = 0.05169880575272146
$ws.Range("D10").Value = 0.4942392263607189
$ws.Range("E10").Value = 0.1408160399458005
$ws.Range("G10").Value = 2.406433429169368
$ws.Range("H10").Value = 1.952681306578626
$ws.Range("I10").Value = 1.925212454088644
$ws.Range("J10").Value = 0.0415877572492267
$ws.Range("K10").Value = 1.192497450114928
$ws.Range("L10").Value = 0.5559591823768528
$ws.Range("N10").Value = 3.043739036098081

$ws.Range("B11").Value = 1.54303060442021
$ws.Range("C11").Value = 0.0553369407261215
$ws.Range("D11").Value = 0.4971111215399304
$ws.Range("E11").Value = 0.1409695529629538
$ws.Range("G11").Value = 2.410980893612106
$ws.Range("H11").Value = 1.951157351613631
$ws.Range("I11").Value = 1.924076409813615
$ws.Range("J11").Value = 0.04153260858088359
$ws.Range("K11").Value = 1.226256644820012
$ws.Range("L11").Value = 0.5604603067621952
$ws.Range("N11").Value = 3.029564251361492

$ws.Range("B12").Value = 1.555693177619219
$ws.Range("C12").Value = 0.05671518812208376
$ws.Range("D12").Value = 0.4982284394370851
$ws.Range("E12").Value = 0.1410344898095595
$ws.Range("G12").Value = 2.412860386200947
$ws.Range("H12").Value = 1.950683130302224
$ws.Range("I12").Value = 1.923754731624619
$ws.Range("J12").Value = 0.041512223373946
$ws.Range("K12").Value = 1.239121529289889
$ws.Range("L12").Value = 0.5621998785414632
$ws.Range("N12").Value = 3.02431666005085

$ws.Range("B13").Value = 1.5529615939368
$ws.Range("C13").Value = 0.05641833233862315
$ws.Range("D13").Value = 0.4979864812870858
$ws.Range("E13").Value = 0.1410202022175078
$ws.Range("G13").Value = 2.412448599718147
$ws.Range("H13").Value = 1.950780688861727
$ws.Range("I13").Value = 1.923819185098097
$ws.Range("J13").Value = 0.04151659155773091
$ws.Range("K13").Value = 1.236347248793948
$ws.Range("L13").Value = 0.5618236713053193
$ws.Range("N13").Value = 3.025441480517493

$ws.Range("B14").Value = 1.544070346052365
$ws.Range("C14").Value = 0.05545031849604243
$ws.Range("D14").Value = 0.4972024473257619
$ws.Range("E14").Value = 0.140974759109433
$ws.Range("G14").Value = 2.411132364463498
$ws.Range("H14").Value = 1.95111627609208
$ws.Range("I14").Value = 1.924047770545087
$ws.Range("J14").Value = 0.04153092150006277
$ws.Range("K14").Value = 1.22731342502513
$ws.Range("L14").Value = 0.5606027195796628
$ws.Range("N14").Value = 3.029130121653864

$ws.Range("B15").Value = 1.538637302952793
$ws.Range("C15").Value = 0.05485745640838502
$ws.Range("D15").Value = 0.4967260814116798
$ws.Range("E15").Value = 0.1409478093986998
$ws.Range("G15").Value = 2.410346640556668
$ws.Range("H15").Value = 1.95133522639253
$ws.Range("I15").Value = 1.924201916646879
$ws.Range("J15").Value = 0.04153976385215863
$ws.Range("K15").Value = 1.221790486763581
$ws.Range("L15").Value = 0.5598594188798103
$ws.Range("N15").Value = 3.031405162395977

$ws.Range("B16").Value = 1.507703932613822
$ws.Range("C16").Value = 0.05146112260614188
$ws.Range("D16").Value = 0.4940557149677147
$ws.Range("E16").Value = 0.1408069609920943
$ws.Range("G16").Value = 2.406158274272173
$ws.Range("H16").Value = 1.952795300818337
$ws.Range("I16").Value = 1.92530188490386
$ws.Range("J16").Value = 0.04159143119316155
$ws.Range("K16").Value = 1.19030256238895
$ws.Range("L16").Value = 0.555669939204023
$ws.Range("N16").Value = 3.044682193099746

$ws.Range("B17").Value = 1.488907252679553
$ws.Range("C17").Value = 0.04937855839844474
$ws.Range("D17").Value = 0.4924706869315969
$ws.Range("E17").Value = 0.1407326998455041
$ws.Range("G17").Value = 2.403869265339097
$ws.Range("H17").Value = 1.953874335522556
$ws.Range("I17").Value = 1.926170007165517
$ws.Range("J17").Value = 0.04162401721710873
$ws.Range("K17").Value = 1.171130400617187
$ws.Range("L17").Value = 0.553162434544106
$ws.Range("N17").Value = 3.053041049347186

$ws.Range("B18").Value = 1.478162187458906
$ws.Range("C18").Value = 0.04818107270261862
$ws.Range("D18").Value = 0.49157858200806
$ws.Range("E18").Value = 0.1406944573728239
$ws.Range("G18").Value = 2.402655726134697
$ws.Range("H18").Value = 1.954562371866018
$ws.Range("I18").Value = 1.9267403889788
$ws.Range("J18").Value = 0.04164308745501621
$ws.Range("K18").Value = 1.160156353948423
$ws.Range("L18").Value = 0.5517432213662801
$ws.Range("N18").Value = 3.057927402596007

$ws.Range("B19").Value = 1.474535493249789
$ws.Range("C19").Value = 0.04777568393267018
$ws.Range("D19").Value = 0.491279892832992
$ws.Range("E19").Value = 0.1406822776578096
$ws.Range("G19").Value = 2.402262538183749
$ws.Range("H19").Value = 1.954806908972301
$ws.Range("I19").Value = 1.926945716560539
$ws.Range("J19").Value = 0.04164960063781109
$ws.Range("K19").Value = 1.15644989073138
$ws.Range("L19").Value = 0.5512666596704463
$ws.Range("N19").Value = 3.059595328724306

$ws.Range("B20").Value = 1.49090133441365
$ws.Range("C20").Value = 0.04960021447796237
$ws.Range("D20").Value = 0.4926373919618641
$ws.Range("E20").Value = 0.1407401425510457
$ws.Range("G20").Value = 2.404102269853098
$ws.Range("H20").Value = 1.953752495089077
$ws.Range("I20").Value = 1.926070239845075
$ws.Range("J20").Value = 0.04162051448636817
$ws.Range("K20").Value = 1.17316579886085
$ws.Range("L20").Value = 0.5534269791101423
$ws.Range("N20").Value = 3.052143103670268

$ws.Range("B21").Value = 1.546679191323506
$ws.Range("C21").Value = 0.05573463208526164
$ws.Range("D21").Value = 0.4974319291815874
$ws.Range("E21").Value = 0.1409879223541957
$ws.Range("G21").Value = 2.411514700689196
$ws.Range("H21").Value = 1.951014914983915
$ws.Range("I21").Value = 1.923977684660429
$ws.Range("J21").Value = 0.0415266989375267
$ws.Range("K21").Value = 1.229964682213847
$ws.Range("L21").Value = 0.5609603910028227
$ws.Range("N21").Value = 3.02804341853232

$ws.Range("B22").Value = 1.58372013331109
$ws.Range("C22").Value = 0.05974713244604857
$ws.Range("D22").Value = 0.5007390574751156
$ws.Range("E22").Value = 0.1411895137790502
$ws.Range("G22").Value = 2.417277011307021
$ws.Range("H22").Value = 1.949825293010718
$ws.Range("I22").Value = 1.923242581680192
$ws.Range("J22").Value = 0.04146828920607915
$ws.Range("K22").Value = 1.26755806532546
$ws.Range("L22").Value = 0.5660884213105959
$ws.Range("N22").Value = 3.012992902689184

$ws.Range("B23").Value = 1.56389714909858
$ws.Range("C23").Value = 0.05760527417668015
$ws.Range("D23").Value = 0.4989581219933683
$ws.Range("E23").Value = 0.1410782999555913
$ws.Range("G23").Value = 2.414117557702525
$ws.Range("H23").Value = 1.950405390752223
$ws.Range("I23").Value = 1.923577059325567
$ws.Range("J23").Value = 0.04149919848653472
$ws.Range("K23").Value = 1.247450692639944
$ws.Range("L23").Value = 0.5633328128316464
$ws.Range("N23").Value = 3.020961571561507

$ws.Range("B24").Value = 1.489999618777631
$ws.Range("C24").Value = 0.04950000436566881
$ws.Range("D24").Value = 0.4925619649664981
$ws.Range("E24").Value = 0.1407367638372534
$ws.Range("G24").Value = 2.40399660940426
$ws.Range("H24").Value = 1.95380736829631
$ws.Range("I24").Value = 1.926115122584555
$ws.Range("J24").Value = 0.04162209702394071
$ws.Range("K24").Value = 1.172245444887125
$ws.Range("L24").Value = 0.5533073087798925
$ws.Range("N24").Value = 3.052548813411839

$ws.Range("B25").Value = 1.412825185954773
$ws.Range("C25").Value = 0.04077040301595503
$ws.Range("D25").Value = 0.4864108443779713
$ws.Range("E25").Value = 0.1405393708336113
$ws.Range("G25").Value = 2.39701860997863
$ws.Range("H25").Value = 1.960054339033746
$ws.Range("I25").Value = 1.931569783396291
$ws.Range("J25").Value = 0.041767196815929
$ws.Range("K25").Value = 1.093162253309117
$ws.Range("L25").Value = 0.5433742628909499
$ws.Range("N25").Value = 3.089624125118185

Write-Host "done"